$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.257.93'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '2.319.96'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'542.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").Value = "'132.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.70%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("E8").Value = '  +2.69%  '

$ws.Range("D9").Value = '2.320.49'
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("E10").Value = '  -0.81%  '

$ws.Range("D11").Value = "'5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.26%  '

$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").Value = "'23.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.32%  '

$ws.Range("D15").Value = '2.735.43'
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").Value = '59.124.25'
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = '2.334.25'
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").Value = "'4.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.83%  '

$ws.Range("D21").Value = "'313.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").Value = "'6.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.42%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").Value = "'62.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '

$ws.Range("D25").Value = "'0.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.25%  '

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").Value = "'7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.50%  '

$ws.Range("E28").Value = '  -1.26%  '

$ws.Range("D29").Value = "'1.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.60%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = "'170.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.21%  '

$ws.Range("B31").Value = 'SuiNetwork'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.21%  '

$ws.Range("D32").Value = '0.0₃0742'
$ws.Range("E32").Value = '  +2.43%  '

$ws.Range("D33").Value = "'5.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.80%  '

$ws.Range("E34").Value = '  +1.41%  '

$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = "'17.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'1.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.84%  '

$ws.Range("E38").Value = '  +0.10%  '

$ws.Range("D39").Value = "'4.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.11%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = "'38.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.00%  '

$ws.Range("D41").Value = "'309.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.77%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.14%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = "'141.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = "'3.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.80%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = "'0.0960"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = "'0.0497"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.58%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'0.558"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = "'18.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.58%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.48%  '

$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = "'11.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D51").Value = "'4.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.10%  '
